$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.647.30"
$ws.Range("E2").Value = "  +0.07%  "
$ws.Range("D3").Value = "1.643.64"
$ws.Range("E3").Value = "  +0.78%  "
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.42"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.00%  "
$ws.Range("E6").Value = "  +1.37%  "
$ws.Range("E8").Value = "  +0.39%  "
$ws.Range("E9").Value = "  +0.82%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.26"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.36%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0841"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.21%  "
$ws.Range("D12").Value = "1.871.89"
$ws.Range("E12").Value = "  +0.66%  "
$ws.Range("D13").Value = "1.637.90"
$ws.Range("E13").Value = "  +2.98%  "
$ws.Range("E14").Value = "  +2.09%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.531"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.68%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.44"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.93%  "
$ws.Range("D17").Value = "26.692.02"
$ws.Range("E17").Value = "  +0.25%  "
$ws.Range("D18").Value = "0.0₃0744"
$ws.Range("E18").Value = "  +0.46%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "216.92"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.41%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.01"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.06%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.34"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.76%  "
$ws.Range("E22").Value = "  +2.33%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.50"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.61%  "
$ws.Range("E24").Value = "  +13.97%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.71"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.25%  "
$ws.Range("E26").Value = "  +0.23%  "
$ws.Range("E27").Value = "  -0.69%  "
$ws.Range("E28").Value = "  +4.92%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.77"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.66%  "
$ws.Range("E30").Value = "  +2.62%  "
$ws.Range("E31").Value = "  +0.90%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.39"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.84%  "
$ws.Range("E33").Value = "  +2.45%  "
$ws.Range("D34").Value = "1.276.30"
$ws.Range("E34").Value = "  +4.16%  "
$ws.Range("E35").Value = "  +2.89%  "
$ws.Range("E36").Value = "  +5.00%  "
$ws.Range("E37").Value = "  +0.58%  "
$ws.Range("E38").Value = "  +6.40%  "
$ws.Range("E39").Value = "  +3.20%  "
$ws.Range("E40").Value = "  +0.16%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.816"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.66%  "
$ws.Range("E42").Value = "  -1.50%  "
$ws.Range("E43").Value = "  +2.31%  "
$ws.Range("D44").Value = "1.782.38"
$ws.Range("E44").Value = "  +0.66%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "91.71"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.23%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "59.81"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +8.52%  "
$ws.Range("E47").Value = "  +1.81%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0517"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.81%  "
$ws.Range("E49").Value = "  +2.33%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0970"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.47%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.407"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.67%  "
